# Running all the suites by making runmode as Y in the excel sheets
$wb = $excel.ActiveWorkbook

# "Test Cases" sheet (the "B suite") holds the TCID / Description / Runmode / Results table.
$ws = $wb.Worksheets.Item("Test Cases")

# Rows 2 and 3 (TestCase_B1, TestCase_B2) currently have Runmode = "N".
# Set them to "Y" so that all suites are run.
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "Y"

$wb.Save()
